$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "(CLIENTE xx/20)" suffixes from the client names in column B
$ws.Range("B5").Value  = "GENIVALDO FREITAS"
$ws.Range("B6").Value  = "PROFESSOR"
$ws.Range("B7").Value  = "JOHTRADER"
$ws.Range("B8").Value  = "SERGIO"
$ws.Range("B9").Value  = "EDSON  BARRETO"
$ws.Range("B10").Value = "WALTER FREIRE"
$ws.Range("B11").Value = "HAREF MACIEL"

# Row 12 was blank (only formatted) - fill it in as a new client entry,
# matching the formatting used by the row above it (row 11)
$ws.Range("B11:E11").Copy()
$ws.Range("B12:E12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B12").Value = "BETO VANN"
$ws.Range("C12").Value = "78344205dc654f4716888e8fb4e3cf10"
$ws.Range("D12").Value = 44827
$ws.Range("E12").Value = 8
